$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update D3 (Segundo Evento) type from "Por defecto" to "Prueba"
$ws.Range("D3").Value = "Prueba"

# Update D4 (Tercer Evento) type from "Por defecto" to "Reunion"
$ws.Range("D4").Value = "Reunion"

# Update the selected/active cell to D4 to match the sheetView selection
$ws.Range("D4").Select()
